$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.117801666259766
$ws.Range("B1").Value = 1.779040694236755
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 1.900094509124756
$ws.Range("E1").Value = 1.115370035171509
